$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(2)
$shp = $hdr.Range.InlineShapes.Item(1)
try { 
  $sr = $shp.Range.ShapeRange
  Write-Host "ShapeRange count:" $sr.Count
  Write-Host "ShapeRange(1).Name:" $sr.Item(1).Name
} catch {
  Write-Host "ERR:" $_.Exception.Message
}
